$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.934.81"
$ws.Range("E2").Value = "  -2.10%  "
$ws.Range("D3").Value = "1.867.48"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4978"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.50%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3807"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08900"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -9.75%  "
$ws.Range("E10").Value = "  -3.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.43"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.305"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.53%  "
$ws.Range("D14").Value = "1.863.29"
$ws.Range("E14").Value = "  -2.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.225"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001097"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.80"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06620"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.98%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.093"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.52%  "
$ws.Range("D23").Value = "27.956.30"
$ws.Range("E23").Value = "  -2.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.284"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.47%  "
$ws.Range("D26").Value = "2.072.96"
$ws.Range("E26").Value = "  -3.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.513"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "157.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("E29").Value = "  -2.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.055"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.12%  "
$ws.Range("E33").Value = "  -2.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.588"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.346"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06534"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02393"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2189"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.289"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.41%  "
$ws.Range("E40").Value = "  -5.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6363"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.91%  "
$ws.Range("E43").Value = "  -3.97%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.92%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5993"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.280"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("E48").Value = "  -2.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.221"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.972"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.53%  "
